# Rename ITEM sheet to EQUIP and update headers/descriptions
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "EQUIP"

$ws1.Range("A1").Value = "equipId"
$ws1.Range("B1").Value = "nameId"
$ws1.Range("C1").Value = "descId"

$ws1.Range("A4").Value = "장비 ID (정본)"
$ws1.Range("B4").Value = "이름 텍스트 ID"
$ws1.Range("C4").Value = "설명 텍스트 ID"

# Add a new CARD sheet after EQUIP
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "CARD"

$ws2.Range("A1").Value = "cardId"
$ws2.Range("B1").Value = "nameId"
$ws2.Range("C1").Value = "descId"

$ws2.Range("A2").Value = "string"
$ws2.Range("B2").Value = "string"
$ws2.Range("C2").Value = "string"

$ws2.Range("A3").Value = "pk"

$ws2.Range("A4").Value = "카드 ID (정본)"
$ws2.Range("B4").Value = "이름 텍스트 ID"
$ws2.Range("C4").Value = "설명 텍스트 ID"
